# Implement init Account & Employee data
# Shifts the existing fullName/nickName/gender/dateOfBirth header row down
# one row and right one column, then fills in the new userName/Account
# columns plus a metadata title row and a phoneNumber column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: insert a blank row above row 1, then a blank column before A ---
$ws.Rows.Item(1).Insert()
$ws.Columns.Item(1).Insert()

# --- Column A: userName ---
$ws.Range("A2").Value = "userName"
$ws.Range("A3").Value = "admin"
$ws.Range("A4").Value = "edwina95"
$ws.Range("A5").Value = "leo.walsh91"
$ws.Range("A6").Value = "kiara.collier"

# --- Column B: fullName (header already there) ---
$ws.Range("B3").Value = "Admin Jack"
$ws.Range("B4").Value = "Edwina"
$ws.Range("B5").Value = "Leo Walsh"
$ws.Range("B6").Value = "Kiara Collier"

# --- Column C: nickName (header already there) ---
$ws.Range("C3").Value = "Jack Jack"
$ws.Range("C4").Value = "Edwina"
$ws.Range("C5").Value = "Leo Walsh"
$ws.Range("C6").Value = "Kiara Collier"

# --- Column D: gender (header already there) ---
$ws.Range("D3").Value = "male"
$ws.Range("D4").Value = "female"
$ws.Range("D5").Value = "male"
$ws.Range("D6").Value = "male"

# --- Column E: dateOfBirth (header already there), stored as text ---
$ws.Range("E3").Value = "'28/03/2000"
$ws.Range("E4").Value = "'29/02/2000"
$ws.Range("E5").Value = "'12/09/2000"

# --- Column F: phoneNumber (new header + numbers) ---
$ws.Range("F2").Value = "phoneNumber"
$ws.Range("F3").Value = 382319487
$ws.Range("F4").Value = 905427980
$ws.Range("F5").Value = 580124839
$ws.Range("F6").Value = 904893167

# --- Row 1: metadata banner ---
$ws.Range("A1").Value = "METADATA: Employee Init Data"
$ws.Range("B1").Value = "userName is from Account "

# --- E6 last (keeps shared-string ordering identical to the authored file) ---
$ws.Range("E6").Value = "'03/12/2000"

# --- Formatting ---
# Header row (2) bold
$ws.Range("A2:F2").Font.Bold = $true

# dateOfBirth column quote-prefixed text, most rows plain, last row date-formatted
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E4").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "mm-dd-yy"

# Title row: vertical-top alignment + taller row
$ws.Range("A1:B1").VerticalAlignment = -4160
$ws.Rows.Item(1).RowHeight = 59.25

# Column widths (engine rounds to a 1/6-character pixel grid, so these are
# the closest achievable values to the authored 29.140625 / 25.140625 / ...)
$ws.Columns.Item(1).ColumnWidth = 28.325
$ws.Columns.Item(2).ColumnWidth = 24.325
$ws.Columns.Item(3).ColumnWidth = 22.0
$ws.Columns.Item(4).ColumnWidth = 19.825
$ws.Columns.Item(5).ColumnWidth = 15.15
$ws.Columns.Item(6).ColumnWidth = 18.65

# Page orientation
$ws.PageSetup.Orientation = 1

# Selection the author ended up on
$ws.Range("F11").Select() | Out-Null
